$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TODOS")
$ws.Range("A2").Value = '18:00'
$ws.Range("B2").Value = '215C_LA PLATA'
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = '🚌'
$ws.Range("A3").Value = '18:04'
$ws.Range("B3").Value = '17_ROMERO'
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = '📅'
$ws.Range("A4").Value = '18:04'
$ws.Range("B4").Value = '23_HERNANDEZ'
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = '🚌'
$ws.Range("A5").Value = '18:08'
$ws.Range("B5").Value = '14_ABASTO'
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = '🚌'
$ws.Range("A6").Value = '18:10'
$ws.Range("B6").Value = '16_SANTA ANA'
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = '🚌'
$ws.Range("A7").Value = '18:16'
$ws.Range("B7").Value = '10_OLMOS'
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = '🚌'
$ws.Range("A8").Value = '18:16'
$ws.Range("B8").Value = '15_ABASTO'
$ws.Range("C8").Value = 18
$ws.Range("D8").Value = '🚌'
$ws.Range("A9").Value = '18:20'
$ws.Range("B9").Value = '16_SANTA ANA'
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = '🚌'
$ws.Range("A10").Value = '18:21'
$ws.Range("B10").Value = '26_HERNANDEZ'
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = '🚌'
$ws.Range("A11").Value = '18:24'
$ws.Range("B11").Value = '14_ABASTO'
$ws.Range("C11").Value = 26
$ws.Range("D11").Value = '🚌'
$ws.Range("A12").Value = '18:27'
$ws.Range("B12").Value = '215C_EL PATO'
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = '🚌'
$ws.Range("A13").Value = '18:31'
$ws.Range("B13").Value = '11X44_ETCHEVERRY'
$ws.Range("C13").Value = 33
$ws.Range("D13").Value = '🚌'
$ws.Range("A14").Value = '18:34'
$ws.Range("B14").Value = '23_HERNANDEZ'
$ws.Range("C14").Value = 36
$ws.Range("D14").Value = '🚌'
$ws.Range("A15").Value = '18:40'
$ws.Range("B15").Value = '15_ABASTO'
$ws.Range("C15").Value = 42
$ws.Range("D15").Value = '🚌'
$ws.Range("A16").Value = '18:47'
$ws.Range("B16").Value = '14X44_ABASTO'
$ws.Range("C16").Value = 49
$ws.Range("D16").Value = '🚌'
$ws.Range("A17").Value = '18:51'
$ws.Range("B17").Value = '215A_LA PLATA'
$ws.Range("C17").Value = 53
$ws.Range("D17").Value = '🚌'
$ws.Range("A18").Value = '18:56'
$ws.Range("B18").Value = '10_OLMOS'
$ws.Range("C18").Value = 58
$ws.Range("D18").Value = '🚌'
$ws.Range("A19").Value = '18:58'
$ws.Range("B19").Value = '215A_EL PATO'
$ws.Range("C19").Value = 60
$ws.Range("D19").Value = '📅'
$ws.Range("A20").Value = '19:04'
$ws.Range("B20").Value = '11_ETCHEVERRY'
$ws.Range("C20").Value = 66
$ws.Range("D20").Value = '📅'
$ws.Range("A21").Value = '19:10'
$ws.Range("B21").Value = '16_P MOR-SANTA ANA'
$ws.Range("C21").Value = 72
$ws.Range("D21").Value = '🚌'
$ws.Range("A22").Value = '19:12'
$ws.Range("B22").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("C22").Value = 74
$ws.Range("D22").Value = '🚌'
$ws.Range("A23").Value = '19:16'
$ws.Range("B23").Value = '27_EL RETIRO'
$ws.Range("C23").Value = 78
$ws.Range("D23").Value = '🚌'
$ws.Range("A24").Value = '19:29'
$ws.Range("B24").Value = '225_GOMEZ'
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = '📅'
$ws.Range("A25").Value = '19:39'
$ws.Range("B25").Value = '215C_EL PATO'
$ws.Range("C25").Value = 101
$ws.Range("D25").Value = '🚌'
$ws.Range("A26").Value = '19:47'
$ws.Range("B26").Value = '26_HERNANDEZ'
$ws.Range("C26").Value = 109
$ws.Range("D26").Value = '🚌'
$ws.Range("A27").Value = '19:49'
$ws.Range("B27").Value = '11X44_ETCHEVERRY'
$ws.Range("C27").Value = 111
$ws.Range("D27").Value = '🚌'
$ws.Range("A28").Value = '19:50'
$ws.Range("B28").Value = '16_P MOR-SANTA ANA'
$ws.Range("C28").Value = 112
$ws.Range("D28").Value = '🚌'
$ws.Range("A29").Value = '19:51'
$ws.Range("B29").Value = '81_EL PELIGRO'
$ws.Range("C29").Value = 113
$ws.Range("D29").Value = '🚌'
$ws.Range("A30").Value = '19:53'
$ws.Range("B30").Value = '215C_LA PLATA'
$ws.Range("C30").Value = 115
$ws.Range("D30").Value = '🚌'

$ws = $wb.Worksheets.Item("215")
$ws.Range("A2").Value = '18:00'
$ws.Range("B2").Value = '215C_LA PLATA'
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = '🚌'
$ws.Range("A3").Value = '18:27'
$ws.Range("B3").Value = '215C_EL PATO'
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = '🚌'
$ws.Range("A4").Value = '18:51'
$ws.Range("B4").Value = '215A_LA PLATA'
$ws.Range("C4").Value = 53
$ws.Range("D4").Value = '🚌'
$ws.Range("A5").Value = '18:58'
$ws.Range("B5").Value = '215A_EL PATO'
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = '📅'
$ws.Range("A6").Value = '19:12'
$ws.Range("B6").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("C6").Value = 74
$ws.Range("D6").Value = '🚌'
$ws.Range("A7").Value = '19:39'
$ws.Range("B7").Value = '215C_EL PATO'
$ws.Range("C7").Value = 101
$ws.Range("D7").Value = '🚌'
$ws.Range("A8").Value = '19:53'
$ws.Range("B8").Value = '215C_LA PLATA'
$ws.Range("C8").Value = 115
$ws.Range("D8").Value = '🚌'

$ws = $wb.Worksheets.Item("COMBINADAS")
$ws.Range("A2").Value = '18:00'
$ws.Range("B2").Value = '215C_LA PLATA'
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = '🚌'
$ws.Range("A3").Value = '18:04'
$ws.Range("B3").Value = '17_ROMERO'
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = '📅'
$ws.Range("A4").Value = '18:04'
$ws.Range("B4").Value = '23_HERNANDEZ'
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = '🚌'
$ws.Range("A5").Value = '18:08'
$ws.Range("B5").Value = '14_ABASTO'
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = '🚌'
$ws.Range("A6").Value = '18:10'
$ws.Range("B6").Value = '16_SANTA ANA'
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = '🚌'
$ws.Range("A7").Value = '18:16'
$ws.Range("B7").Value = '10_OLMOS'
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = '🚌'
$ws.Range("A8").Value = '18:16'
$ws.Range("B8").Value = '15_ABASTO'
$ws.Range("C8").Value = 18
$ws.Range("D8").Value = '🚌'
$ws.Range("A9").Value = '18:20'
$ws.Range("B9").Value = '16_SANTA ANA'
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = '🚌'
$ws.Range("A10").Value = '18:21'
$ws.Range("B10").Value = '26_HERNANDEZ'
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = '🚌'
$ws.Range("A11").Value = '18:24'
$ws.Range("B11").Value = '14_ABASTO'
$ws.Range("C11").Value = 26
$ws.Range("D11").Value = '🚌'
$ws.Range("A12").Value = '18:27'
$ws.Range("B12").Value = '215C_EL PATO'
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = '🚌'
$ws.Range("A13").Value = '18:31'
$ws.Range("B13").Value = '11X44_ETCHEVERRY'
$ws.Range("C13").Value = 33
$ws.Range("D13").Value = '🚌'
$ws.Range("A14").Value = '18:34'
$ws.Range("B14").Value = '23_HERNANDEZ'
$ws.Range("C14").Value = 36
$ws.Range("D14").Value = '🚌'
$ws.Range("A15").Value = '18:40'
$ws.Range("B15").Value = '15_ABASTO'
$ws.Range("C15").Value = 42
$ws.Range("D15").Value = '🚌'
$ws.Range("A16").Value = '18:47'
$ws.Range("B16").Value = '14X44_ABASTO'
$ws.Range("C16").Value = 49
$ws.Range("D16").Value = '🚌'
$ws.Range("A17").Value = '18:51'
$ws.Range("B17").Value = '215A_LA PLATA'
$ws.Range("C17").Value = 53
$ws.Range("D17").Value = '🚌'
$ws.Range("A18").Value = '18:56'
$ws.Range("B18").Value = '10_OLMOS'
$ws.Range("C18").Value = 58
$ws.Range("D18").Value = '🚌'
$ws.Range("A19").Value = '18:58'
$ws.Range("B19").Value = '215A_EL PATO'
$ws.Range("C19").Value = 60
$ws.Range("D19").Value = '📅'
$ws.Range("A20").Value = '19:04'
$ws.Range("B20").Value = '11_ETCHEVERRY'
$ws.Range("C20").Value = 66
$ws.Range("D20").Value = '📅'
$ws.Range("A21").Value = '19:10'
$ws.Range("B21").Value = '16_P MOR-SANTA ANA'
$ws.Range("C21").Value = 72
$ws.Range("D21").Value = '🚌'
$ws.Range("A22").Value = '19:12'
$ws.Range("B22").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("C22").Value = 74
$ws.Range("D22").Value = '🚌'
$ws.Range("A23").Value = '19:16'
$ws.Range("B23").Value = '27_EL RETIRO'
$ws.Range("C23").Value = 78
$ws.Range("D23").Value = '🚌'
$ws.Range("A24").Value = '19:29'
$ws.Range("B24").Value = '225_GOMEZ'
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = '📅'
$ws.Range("A25").Value = '19:39'
$ws.Range("B25").Value = '215C_EL PATO'
$ws.Range("C25").Value = 101
$ws.Range("D25").Value = '🚌'
$ws.Range("A26").Value = '19:47'
$ws.Range("B26").Value = '26_HERNANDEZ'
$ws.Range("C26").Value = 109
$ws.Range("D26").Value = '🚌'
$ws.Range("A27").Value = '19:49'
$ws.Range("B27").Value = '11X44_ETCHEVERRY'
$ws.Range("C27").Value = 111
$ws.Range("D27").Value = '🚌'
$ws.Range("A28").Value = '19:50'
$ws.Range("B28").Value = '16_P MOR-SANTA ANA'
$ws.Range("C28").Value = 112
$ws.Range("D28").Value = '🚌'
$ws.Range("A29").Value = '19:51'
$ws.Range("B29").Value = '81_EL PELIGRO'
$ws.Range("C29").Value = 113
$ws.Range("D29").Value = '🚌'
$ws.Range("A30").Value = '19:53'
$ws.Range("B30").Value = '215C_LA PLATA'
$ws.Range("C30").Value = 115
$ws.Range("D30").Value = '🚌'

